$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New DAP (daily availability plan) file for 2025-01-13 ---
# Only 14 hourly rows this time (vs. 17 previously), and stations S3/eS3 report 0

# Date column for all 14 data rows (A2:A15); quoted so Excel keeps it as literal text,
# matching the inline string used in the source DAP export (not a date serial).
$ws.Range("A2:A15").Value = "'2025-01-13"

# Hour / S1-S8 / Total / eS1-eS8 / eTotal for each row, in column order B..T
$data = @(
    @(1, 20163, 0, 0, 12898, 0, 0, 0, 0, 33061, 20165.0163, 0, 0, 12899.2898, 0, 0, 0, 0, 33064.3061),
    @(2, 19359, 0, 0, 12469, 0, 0, 0, 0, 31828, 19360.9359, 0, 0, 12470.2469, 0, 0, 0, 0, 31831.1828),
    @(3, 18549, 0, 0, 12025, 0, 0, 0, 0, 30574, 18550.8549, 0, 0, 12026.2025, 0, 0, 0, 0, 30577.0574),
    @(4, 17841, 0, 0, 11802, 0, 0, 0, 0, 29643, 17842.7841, 0, 0, 11803.1802, 0, 0, 0, 0, 29645.9643),
    @(5, 18016, 0, 0, 11832, 0, 0, 0, 0, 29848, 18017.8016, 0, 0, 11833.1832, 0, 0, 0, 0, 29850.9848),
    @(6, 18546, 0, 0, 12250, 0, 0, 0, 0, 30796, 18547.8546, 0, 0, 12251.225, 0, 0, 0, 0, 30799.0796),
    @(7, 19203, 0, 0, 12410, 0, 0, 0, 0, 31613, 19204.9203, 0, 0, 12411.241, 0, 0, 0, 0, 31616.1613),
    @(8, 21420, 0, 0, 13519, 0, 0, 0, 0, 34939, 21422.142, 0, 0, 13520.3519, 0, 0, 0, 0, 34942.4939),
    @(9, 26095, 0, 0, 15839, 0, 0, 0, 0, 41934, 26097.6095, 0, 0, 15840.5839, 0, 0, 0, 0, 41938.1934),
    @(10, 28036, 0, 0, 15116, 0, 0, 0, 0, 43152, 28038.8036, 0, 0, 15117.5116, 0, 0, 0, 0, 43156.3152),
    @(11, 28526, 0, 0, 15853, 0, 0, 0, 0, 44379, 28528.8526, 0, 0, 15854.5853, 0, 0, 0, 0, 44383.4379),
    @(12, 29282, 0, 0, 15881, 0, 0, 0, 0, 45163, 29284.9282, 0, 0, 15882.5881, 0, 0, 0, 0, 45167.5163),
    @(13, 28956, 0, 0, 16391, 0, 0, 0, 0, 45347, 28958.8956, 0, 0, 16392.6391, 0, 0, 0, 0, 45351.5347),
    @(14, 30957, 0, 0, 16335, 0, 0, 0, 0, 47292, 30960.0957, 0, 0, 16336.6335, 0, 0, 0, 0, 47296.7292)
)

$r = 2
foreach ($row in $data) {
    $c = 2
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Drop rows 16-18, which only existed in the previous 17-hour DAP file
$ws.Range("A16:T18").Delete()
